$d = $word.ActiveDocument

$replacements = @(
    @{old="68×62=4216"; new="23×32=736"},
    @{old="12×94=1128"; new="41×28=1148"},
    @{old="43×23=989";  new="23×83=1909"},
    @{old="19×40=760";  new="41×76=3116"},
    @{old="60×16=960";  new="25×21=525"},
    @{old="74×37=2738"; new="68×14=952"},
    @{old="34×79=2686"; new="19×97=1843"},
    @{old="48×85=4080"; new="77×95=7315"},
    @{old="75×41=3075"; new="49×15=735"},
    @{old="31×97=3007"; new="75×74=5550"},
    @{old="43×75=3225"; new="59×65=3835"},
    @{old="90×73=6570"; new="42×75=3150"},
    @{old="77×79=6083"; new="53×68=3604"},
    @{old="68×55=3740"; new="72×94=6768"},
    @{old="81×38=3078"; new="98×48=4704"},
    @{old="89×91=8099"; new="99×48=4752"},
    @{old="63×15=945";  new="50×56=2800"},
    @{old="87×87=7569"; new="49×67=3283"},
    @{old="89×61=5429"; new="23×24=552"},
    @{old="94×70=6580"; new="64×25=1600"},
    @{old="79×23=1817"; new="77×32=2464"},
    @{old="64×58=3712"; new="71×78=5538"},
    @{old="74×45=3330"; new="65×96=6240"},
    @{old="67×88=5896"; new="59×27=1593"},
    @{old="14×66=924";  new="16×91=1456"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
